$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRef, $value)
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = $origStyle
}

# Row 2
Set-TextValue "D2" '97.597.90'
Set-TextValue "E2" '  +2.98%  '

# Row 3
Set-TextValue "D3" '3.600.07'
Set-TextValue "E3" '  +1.22%  '

# Row 4
Set-TextValue "E4" '  -0.06%  '

# Row 5
Set-TextValue "D5" '242.56'
Set-TextValue "E5" '  +2.93%  '

# Row 6
Set-TextValue "E6" '  +17.98%  '

# Row 7
Set-TextValue "D7" '654.35'
Set-TextValue "E7" '  +0.20%  '

# Row 8
Set-TextValue "D8" '0.436'
Set-TextValue "E8" '  +10.03%  '

# Row 9
Set-TextValue "E9" '  -0.10%  '

# Row 10
Set-TextValue "D10" '1.05'
Set-TextValue "E10" '  +5.30%  '

# Row 11
Set-TextValue "D11" '3.599.69'
Set-TextValue "E11" '  +1.25%  '

# Row 12
Set-TextValue "D12" '44.49'
Set-TextValue "E12" '  +5.69%  '

# Row 13
Set-TextValue "E13" '  +1.13%  '

# Row 14
Set-TextValue "D14" '6.48'
Set-TextValue "E14" '  +0.84%  '

# Row 15
Set-TextValue "D15" '4.267.34'
Set-TextValue "E15" '  +1.16%  '

# Row 16
Set-TextValue "D16" '97.192.92'
Set-TextValue "E16" '  +2.62%  '

# Row 17
Set-TextValue "E17" '  +4.19%  '

# Row 18
Set-TextValue "D18" '8.67'
Set-TextValue "E18" '  +2.19%  '

# Row 19
Set-TextValue "D19" '3.595.71'
Set-TextValue "E19" '  +1.11%  '

# Row 20
Set-TextValue "D20" '12.60'
Set-TextValue "E20" '  -0.51%  '

# Row 21
Set-TextValue "D21" '18.18'
Set-TextValue "E21" '  +2.82%  '

# Row 22
Set-TextValue "E22" '  +11.20%  '

# Row 23
Set-TextValue "E23" '  +1.98%  '

# Row 24
Set-TextValue "D24" '518.91'
Set-TextValue "E24" '  +2.51%  '

# Row 25
Set-TextValue "D25" '0.0000211'
Set-TextValue "E25" '  +8.66%  '

# Row 26
Set-TextValue "E26" '  +2.48%  '

# Row 27
Set-TextValue "D27" '102.31'
Set-TextValue "E27" '  +8.05%  '

# Row 28
Set-TextValue "D28" '13.15'
Set-TextValue "E28" '  +5.09%  '

# Row 29
Set-TextValue "D29" '3.792.04'
Set-TextValue "E29" '  +1.21%  '

# Row 30
Set-TextValue "D30" '0.169'
Set-TextValue "E30" '  +18.38%  '

# Row 31
Set-TextValue "D31" '12.09'
Set-TextValue "E31" '  +5.85%  '

# Row 32
Set-TextValue "E32" '  -0.84%  '

# Row 33
Set-TextValue "D33" '1.00'
Set-TextValue "E33" '  +0.19%  '

# Row 34
Set-TextValue "E34" '  +6.38%  '

# Row 35
Set-TextValue "D35" '0.999'
Set-TextValue "E35" '  -0.02%  '

# Row 36
Set-TextValue "E36" '  +1.20%  '

# Row 37
Set-TextValue "B37" 'Bittensor'
Set-TextValue "C37" 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue "D37" '620.47'
Set-TextValue "E37" '  +7.16%  '

# Row 38
Set-TextValue "B38" 'PolygonEcosystemToken'
Set-TextValue "C38" 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
Set-TextValue "D38" '0.574'
Set-TextValue "E38" '  +3.91%  '

# Row 39
Set-TextValue "D39" '8.76'
Set-TextValue "E39" '  +4.07%  '

# Row 40
Set-TextValue "E40" '  -2.25%  '

# Row 41
Set-TextValue "E41" '  +3.58%  '

# Row 42
Set-TextValue "E42" '  +8.21%  '

# Row 43
Set-TextValue "D43" '0.933'
Set-TextValue "E43" '  +3.53%  '

# Row 44
Set-TextValue "E44" '  -0.02%  '

# Row 45
Set-TextValue "D45" '6.04'
Set-TextValue "E45" '  +5.61%  '

# Row 46
Set-TextValue "E46" '  +7.63%  '

# Row 47
Set-TextValue "D47" '0.433'
Set-TextValue "E47" '  +43.19%  '

# Row 48
Set-TextValue "D48" '2.32'
Set-TextValue "E48" '  +2.19%  '

# Row 49
Set-TextValue "D49" '23.66'
Set-TextValue "E49" '  +1.24%  '

# Row 50
Set-TextValue "D50" '8.58'
Set-TextValue "E50" '  +5.89%  '

# Row 51
Set-TextValue "D51" '3.32'
Set-TextValue "E51" '  +9.07%  '
